# Apply the "include_RPS" scenario-input update to scenarios.xlsx
$wb = $excel.ActiveWorkbook

# --- SolverSettings sheet: add new row 10 with include_RPS / N values ---
$wsSolver = $wb.Worksheets.Item("SolverSettings")
$wsSolver.Range("A10").Value = "include_RPS"
$wsSolver.Range("B10:G10").Value = "N"

# Select H10 on SolverSettings, and make it the active/visible sheet (tabSelected)
$wsSolver.Activate() | Out-Null
$wsSolver.Range("H10").Select() | Out-Null
